$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the |S*|/n column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the new summary labels / values: bold 12pt font, vertical-center alignment, row height 15.6
$summaryRange = $ws.Range("A14:B17")
$summaryRange.Font.Bold = $true
$summaryRange.Font.Size = 12
$summaryRange.VerticalAlignment = -4108
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Selection as left by the editing session
$ws.Range("A14:B17").Select()

# Page setup used when the workbook was printed/saved
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
